$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain text (matches source data format)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.442.86"

$ws.Range("D3").Value = "2.100.88"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "334.92"
$ws.Range("E5").Value = "  +1.47%  "

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "0.5224"
$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("D8").Value = "0.4562"
$ws.Range("E8").Value = "  +3.91%  "

$ws.Range("D9").Value = "56.29"
$ws.Range("E9").Value = "  +12.85%  "

$ws.Range("D10").Value = "0.08934"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("D11").Value = "1.178"
$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("E12").Value = "  -2.55%  "

$ws.Range("D13").Value = "2.089.46"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").Value = "6.837"
$ws.Range("E14").Value = "  +1.38%  "

$ws.Range("D15").Value = "8.049"
$ws.Range("E15").Value = "  +3.72%  "

$ws.Range("D16").Value = "97.32"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("D17").Value = "0.00001154"
$ws.Range("E17").Value = "  +2.10%  "

$ws.Range("D19").Value = "0.06634"
$ws.Range("E19").Value = "  -0.18%  "

$ws.Range("D20").Value = "19.20"
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").Value = "6.302"
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("D23").Value = "30.494.66"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").Value = "12.42"
$ws.Range("E24").Value = "  +1.15%  "

$ws.Range("D25").Value = "2.356"
$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("D26").Value = "2.343.44"
$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("E27").Value = "  -1.15%  "

$ws.Range("D28").Value = "163.08"
$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("E29").Value = "  -4.27%  "

$ws.Range("D30").Value = "133.32"
$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("D31").Value = "1.214"
$ws.Range("E31").Value = "  -0.77%  "

$ws.Range("D32").Value = "0.1069"

$ws.Range("D33").Value = "1.664"
$ws.Range("E33").Value = "  -1.32%  "

$ws.Range("D34").Value = "6.356"
$ws.Range("E34").Value = "  +1.89%  "

$ws.Range("D35").Value = "3.945"
$ws.Range("E35").Value = "  +1.25%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "5.963"
$ws.Range("E36").Value = "  +8.10%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "10.29"
$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("D38").Value = "0.02584"
$ws.Range("E38").Value = "  -0.23%  "

$ws.Range("D39").Value = "0.06854"
$ws.Range("E39").Value = "  +1.54%  "

$ws.Range("D40").Value = "0.2335"
$ws.Range("E40").Value = "  +2.31%  "

$ws.Range("D41").Value = "12.66"
$ws.Range("E41").Value = "  -0.59%  "

$ws.Range("D42").Value = "0.6878"
$ws.Range("E42").Value = "  -0.72%  "

$ws.Range("D43").Value = "1.247"
$ws.Range("E43").Value = "  -2.02%  "

$ws.Range("E44").Value = "  +4.57%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.6408"
$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "13.98"
$ws.Range("E46").Value = "  -0.90%  "

$ws.Range("D47").Value = "3.660"
$ws.Range("E47").Value = "  +0.60%  "

$ws.Range("E48").Value = "  -0.32%  "

$ws.Range("E49").Value = "  +14.77%  "

$ws.Range("B50").Value = "WEMIXTOKEN"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "1.203"
$ws.Range("E50").Value = "  -1.43%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "83.03"
$ws.Range("E51").Value = "  +0.24%  "

Write-Output "done"
